# Tester function for Circuit Three: fix the wrong output value shown for
# the "1 0 0 1 0 1" row in the Input/Output table on Sheet2.
#
# Row 7 (the "wrong output" row) had E7 = "1 0 0 1 0 1"; the corrected
# value is "0 1 1 0 0 1" (matching F7 / Circuit4's value for that input).
# Updating the cell drops the now-unused shared string and the active
# selection in the sheet is left on E8, matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "0 1 1 0 0 1"

$ws.Range("E8").Select() | Out-Null
